$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text in A1 ("Scenario: Test Suite" -> "Scenario: Scenario")
$ws.Range("A1").Value = "Scenario: Scenario"

# The sheet lost data: rows 4 and 5 (each A:E, mirroring row 3) are missing.
# Restore them by duplicating row 3 (which already holds the blank-cell
# pattern used by the rest of the table) down into rows 4 and 5, so the
# new cells pick up the same content/style as row 3 instead of being
# built up from scratch.
$ws.Range("A3:E3").Copy($ws.Range("A4:E4"))
$ws.Range("A3:E3").Copy($ws.Range("A5:E5"))
